$wb = $excel.ActiveWorkbook

# Rename header labels on Sim1 (Coal1 -> GasPlant1, Gas1 -> CoalPlant1, Consumer1 -> HydroPlant1)
$ws1 = $wb.Worksheets.Item("Sim1")
$ws1.Range("B1").Value = "GasPlant1"
$ws1.Range("C1").Value = "CoalPlant1"
$ws1.Range("D1").Value = "HydroPlant1"

# Add new column E header
$ws1.Range("E1").Value = "ZeroBidders"

# Remove the Sim2 sheet entirely
$ws2 = $wb.Worksheets.Item("Sim2")
$ws2.Delete()

# Make Sim1 the active/selected sheet and select cell E2
$ws1.Select()
$ws1.Range("E2").Select()
